$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 397
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H19").Value = 1733
$ws.Range("I19").Value = 1300
$ws.Range("K19").Value = 1300
$ws.Range("M19").Value = -1125
$ws.Range("H44").Value = 3092
$ws.Range("J44").Value = 3092
$ws.Range("L44").Value = 3092
$ws.Range("N44").Value = -4016
$ws.Range("H113").Value = 9242.375
$ws.Range("J113").Value = 10039.8
$ws.Range("L113").Value = 10039.8
$ws.Range("N113").Value = -16547.8
$ws.Range("H132").Value = 26177.5
$ws.Range("I132").Value = 4148.5625
$ws.Range("K132").Value = 12445.6875
$ws.Range("M132").Value = -9915.6875
$ws.Range("H138").Value = 3593.9465
$ws.Range("I138").Value = 1610.2593
$ws.Range("J138").Value = 5440.8276
$ws.Range("K138").Value = 4830.7779
$ws.Range("L138").Value = 16322.4828
$ws.Range("M138").Value = 309.2221
$ws.Range("N138").Value = -26602.4828
$ws.Range("H141").Value = 3996.0967
$ws.Range("I141").Value = 1958.5186
$ws.Range("J141").Value = 17749.75
$ws.Range("K141").Value = 5875.5558
$ws.Range("L141").Value = 53249.25
$ws.Range("M141").Value = -695.5558000000001
$ws.Range("N141").Value = -63609.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 45571.285
$ws.Range("I2").Value = 52499.832
$ws.Range("K2").Value = 52499.832
$ws.Range("M2").Value = -52386.832
$ws.Range("H74").Value = 2272.6785
$ws.Range("I74").Value = 2101.4
$ws.Range("K74").Value = 2101.4
$ws.Range("M74").Value = -1227.4
$ws.Range("H77").Value = 2272.6785
$ws.Range("I77").Value = 2101.4
$ws.Range("K77").Value = 10507
$ws.Range("M77").Value = -6139
$ws.Range("H116").Value = 45571.285
$ws.Range("I116").Value = 52499.832
$ws.Range("K116").Value = 52499.832
$ws.Range("M116").Value = -50205.832
$ws.Range("H122").Value = 2245.4167
$ws.Range("I122").Value = 2245.4167
$ws.Range("K122").Value = 6736.250100000001
$ws.Range("M122").Value = -4286.250100000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 45571.285
$ws.Range("I3").Value = 52499.832
$ws.Range("K3").Value = 52499.832
$ws.Range("M3").Value = -52385.832
$ws.Range("H134").Value = 1295.8718
$ws.Range("I134").Value = 1290.5
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 3871.5
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -1336.5
$ws.Range("N134").Value = -9570

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1622.898
$ws.Range("J31").Value = 2196
$ws.Range("L31").Value = 2196
$ws.Range("N31").Value = -2786
$ws.Range("H34").Value = 1622.898
$ws.Range("J34").Value = 2196
$ws.Range("L34").Value = 2196
$ws.Range("N34").Value = -2600
$ws.Range("H99").Value = 28055882
$ws.Range("I99").Value = 4069507.2
$ws.Range("K99").Value = 4069507.2
$ws.Range("M99").Value = -4068009.2
$ws.Range("H122").Value = 381509.8
$ws.Range("I122").Value = 639429.2
$ws.Range("J122").Value = 6354.364
$ws.Range("K122").Value = 1918287.6
$ws.Range("L122").Value = 19063.092
$ws.Range("M122").Value = -1915837.6
$ws.Range("N122").Value = -23963.092
$ws.Range("H126").Value = 28055882
$ws.Range("I126").Value = 4069507.2
$ws.Range("K126").Value = 12208521.6
$ws.Range("M126").Value = -12206051.6
$ws.Range("H134").Value = 2899.7812
$ws.Range("I134").Value = 1995.5
$ws.Range("J134").Value = 16464
$ws.Range("K134").Value = 5986.5
$ws.Range("L134").Value = 49392
$ws.Range("M134").Value = -3451.5
$ws.Range("N134").Value = -54462

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 37222544
$ws.Range("I4").Value = 41875236
$ws.Range("K4").Value = 125625708
$ws.Range("M4").Value = -125625596
$ws.Range("H108").Value = 1533.3334
$ws.Range("I108").Value = 1533.3334
$ws.Range("K108").Value = 4600.0002
$ws.Range("M108").Value = -1720.0002
$ws.Range("H109").Value = 362.81818
$ws.Range("I109").Value = 362.81818
$ws.Range("K109").Value = 1088.45454
$ws.Range("M109").Value = -48.45453999999995
$ws.Range("H125").Value = 6272.5
$ws.Range("I125").Value = 5545
$ws.Range("K125").Value = 16635
$ws.Range("M125").Value = -11715
$ws.Range("H131").Value = 2908.9614
$ws.Range("I131").Value = 1791.4667
$ws.Range("J131").Value = 4432.8184
$ws.Range("K131").Value = 5374.4001
$ws.Range("L131").Value = 13298.4552
$ws.Range("M131").Value = -334.4000999999998
$ws.Range("N131").Value = -23378.4552

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5230.375
$ws.Range("I122").Value = 4153.4546
$ws.Range("K122").Value = 12460.3638
$ws.Range("M122").Value = -10010.3638
$ws.Range("H126").Value = 6300.6665
$ws.Range("I126").Value = 5850
$ws.Range("K126").Value = 17550
$ws.Range("M126").Value = -15080
$ws.Range("H132").Value = 3618.6428
$ws.Range("I132").Value = 3532.818
$ws.Range("K132").Value = 10598.454
$ws.Range("M132").Value = -8068.454000000002

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1032.5834
$ws.Range("I16").Value = 741.625
$ws.Range("J16").Value = 1614.5
$ws.Range("K16").Value = 741.625
$ws.Range("L16").Value = 1614.5
$ws.Range("M16").Value = -571.625
$ws.Range("N16").Value = -1954.5
$ws.Range("H68").Value = 2848.8948
$ws.Range("I68").Value = 2888.25
$ws.Range("K68").Value = 2888.25
$ws.Range("M68").Value = -2139.25
$ws.Range("H71").Value = 2848.8948
$ws.Range("I71").Value = 2888.25
$ws.Range("K71").Value = 14441.25
$ws.Range("M71").Value = -10697.25
$ws.Range("H122").Value = 6986.6665
$ws.Range("I122").Value = 2650
$ws.Range("J122").Value = 7653.846
$ws.Range("K122").Value = 7950
$ws.Range("L122").Value = 22961.538
$ws.Range("M122").Value = -5500
$ws.Range("N122").Value = -27861.538

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1377.0667
$ws.Range("I113").Value = 1222.7778
$ws.Range("J113").Value = 1608.5
$ws.Range("K113").Value = 3668.3334
$ws.Range("L113").Value = 4825.5
$ws.Range("M113").Value = -1498.3334
$ws.Range("N113").Value = -9165.5
$ws.Range("H122").Value = 2158.5
$ws.Range("I122").Value = 2003
$ws.Range("K122").Value = 6009
$ws.Range("M122").Value = -3559
$ws.Range("H126").Value = 1336.76
$ws.Range("I126").Value = 1201.7894
$ws.Range("K126").Value = 3605.3682
$ws.Range("M126").Value = -1135.3682
